$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the existing "_GoBack" bookmark. In the source document it wraps
#    the "valorCredito (ExtensoValorCredito), " text further down the
#    contract; the edit relocates it to surround the word "BMP" in the
#    "CEDENTE:" clause near the top of the document.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Shorten "BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A., " down to
#    "BMP SOCIEDADE DE CRÉDITO DIRETO S.A., " right after "CEDENTE: ".
#
#    The two runs that follow in the same paragraph ("instituição
#    financeira...Estatuto Social" and ";") happen to share identical
#    run formatting, so rewriting any text earlier in the paragraph would
#    otherwise cause them to be silently recombined into a single run.
#    A throw-away bookmark dropped right between them keeps them from
#    being treated as adjacent/mergeable while we edit, and is removed
#    again once the text edit is done.
# ---------------------------------------------------------------------------
$rGuard = $d.Content
$null = $rGuard.Find.Execute("Estatuto Social;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$guardPoint = $d.Range($rGuard.End - 1, $rGuard.End - 1)
$d.Bookmarks.Add("zzTempGuard", $guardPoint)

$r1 = $d.Content
$null = $r1.Find.Execute(
    "BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A., ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "BMP SOCIEDADE DE CRÉDITO DIRETO S.A., ", 2)

$d.Bookmarks("zzTempGuard").Delete()

# ---------------------------------------------------------------------------
# 3) Re-insert "_GoBack" around the "BMP" that immediately follows
#    "CEDENTE: " so it ends up on its own run, exactly as in the edited
#    document.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute(
    "CEDENTE: BMP", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
$bmpRange = $d.Range($r2.End - 3, $r2.End)
$d.Bookmarks.Add("_GoBack", $bmpRange)

# ---------------------------------------------------------------------------
# 4) Apply the same shortening further down, in the signature block:
#    "BMP MONEY PLUS SOCIEDADE DE CRÉDITO" -> "BMP SOCIEDADE DE CRÉDITO".
# ---------------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute(
    "BMP MONEY PLUS SOCIEDADE DE CRÉDITO", $true, $false, $false, $false,
    $false, $true, 1, $false, "BMP SOCIEDADE DE CRÉDITO", 2)

# ---------------------------------------------------------------------------
# 5) Update the cached result of the footer's PAGE field from "6" to "2".
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $footer = $sec.Footers.Item(1)
    foreach ($f in $footer.Range.Fields) {
        if ($f.Code.Text -match "PAGE") {
            $null = $f.Result.Find.Execute(
                "6", $true, $false, $false, $false, $false, $true, 1, $false,
                "2", 2)
        }
    }
}
